$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as TEXT
# (so numeric/date-looking strings like phone numbers, years with leading
# zeros, and ISO dates keep their literal text instead of being silently
# reinterpreted as a number/date). We reset the style back to "Normal"
# afterwards so we don't leave a stray text-format style applied to the
# cell (matches the original file, which has no custom cell styles).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# The source form export leaves several fields genuinely empty, but still
# present as an (empty-string) cell rather than a truly blank cell -
# Range.Value can't directly store a bare "" (Excel/COM treats that as a
# clear-cell request), so give those cells the formula "="""" ", whose
# cached/displayed result is the same empty string.
function Set-EmptyText($range) {
    $range.Formula = "=""""" 
}

# F39 currently holds a formula "=08420880979" whose cached result is the
# text "08420880979". The edit drops the formula and keeps just the
# literal text value.
Set-TextValue $ws.Range("F39") "08420880979"

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N")

# New rows 40-45 appended at the bottom of the sheet.
$rows = @()
$rows += ,@("BSS/3fca2c65e357","AGNIVA","BHATTACHARJEE","nei","2020","08420880979","bhattacharjee.agniva.jobs@gmail.com","2002-01-21","cdvf","","","","","dvfgbdbfd")
$rows += ,@("BSS/90242099a5fc","AGNIVA","BHATTACHARJEE","dfsdvsca","2021","08420880979","bhattacharjee.agniva.jobs@gmail.com","2002-01-21","sdvsfdvd","","","","","sdvsfbgdfsdv")
$rows += ,@("BSS/29a5a21e81a8","AGNIVA","BHATTACHARJEE","dfsdvsca","2021","08420880979","bhattacharjee.agniva.jobs@gmail.com","2002-01-21","sdvsfdvd","","","","","sdvsfbgdfsdv")
$rows += ,@("BSS/8e5c45e9e1d2","AGNIVA","BHATTACHARJEE","dfsdvsca","2021","08420880979","bhattacharjee.agniva.jobs@gmail.com","2002-01-21","sdvsfdvd","","","","","sdvsfbgdfsdv")
$rows += ,@("BSS/16fe73826fad","AGNIVA","BHATTACHARJEE","dfsdvsca","2021","08420880979","bhattacharjee.agniva.jobs@gmail.com","2002-01-21","sdvsfdvd","","","","","sdvsfbgdfsdv")
$rows += ,@("BSS/9b6efdce74cd","AGNIVA","BHATTACHARJEE","biye hobe na er","2020","08420880979","bhattacharjee.agniva.jobs@gmail.com","2002-01-21","IT","Google","","","","dsfdvsdscabg")

$startRow = 40
for ($ri = 0; $ri -lt $rows.Length; $ri++) {
    $rowvals = $rows[$ri]
    $r = $startRow + $ri
    for ($ci = 0; $ci -lt $cols.Length; $ci++) {
        $val = $rowvals[$ci]
        $addr = $cols[$ci] + $r
        if ($val -ne "") {
            Set-TextValue $ws.Range($addr) $val
        } else {
            Set-EmptyText $ws.Range($addr)
        }
    }
}

# Row 45's Phone cell (F45) is stored as a formula "=08420880979" (mirrors
# what F39 used to be before this edit), whose cached text is the same
# phone number.
$ws.Range("F45").Formula = "=08420880979"
